# Femacal de La Calera - Sandia : insert a new weekly price record.
#
# The source data row previously at sheet row 178 (and all rows below it,
# down to row 277) need to shift down by one row to make room for a new
# record. The new record is inserted as the new row 178.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 178:277 down to 179:278 by inserting a new blank row at 178.
$ws.Range("A178:R178").Insert()

# Populate the newly inserted row 178 with the new observation.
$ws.Range("A178").Value() = 3
$ws.Range("B178").Value() = "Femacal de La Calera"
$ws.Range("C178").Value() = "Coquimbo"
$ws.Range("D178").Value() = 44518
$ws.Range("E178").Value() = 5
$ws.Range("F178").Value() = 100112028
$ws.Range("G178").Value() = "Sandia"
$ws.Range("H178").Value() = "Sin especificar"
$ws.Range("I178").Value() = "Primera"
$ws.Range("J178").Value() = 260
$ws.Range("K178").Value() = 700
$ws.Range("L178").Value() = 700
$ws.Range("M178").Value() = 700
$ws.Range("N178").Value() = "$/kilo (volumen en unidades)"
$ws.Range("O178").Value() = "Perú"
$ws.Range("P178").Value() = 700
$ws.Range("Q178").Value() = 1
$ws.Range("R178").Value() = "Hortaliza"

# Make sure the style applied to the date column (D) on the new row matches
# the date style used throughout the rest of the column.
$ws.Range("D178").NumberFormat = $ws.Range("D179").NumberFormat
